$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.926.01"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "3.029.87"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Formula = "'594.46"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Formula = "'153.39"
$ws.Range("E6").Value = "  +8.06%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.023.39"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +17.39%  "
$ws.Range("D11").Formula = "'0.150"
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").Formula = "'35.52"
$ws.Range("E14").Value = "  +4.63%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "3.537.73"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("E17").Value = "  +3.87%  "
$ws.Range("D18").Value = "62.869.19"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").Value = "3.029.71"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Formula = "'452.28"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Formula = "'14.28"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Formula = "'0.697"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("E23").Value = "  +4.12%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Formula = "'83.17"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Formula = "'11.41"
$ws.Range("E25").Value = "  +10.93%  "
$ws.Range("D26").Formula = "'2.32"
$ws.Range("E26").Value = "  +7.96%  "
$ws.Range("D27").Formula = "'12.42"
$ws.Range("E27").Value = "  +4.65%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Formula = "'2.30"
$ws.Range("E29").Value = "  +13.57%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Formula = "'7.53"
$ws.Range("E30").Value = "  +7.35%  "
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Formula = "'27.58"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").Value = "0.0₃0861"
$ws.Range("E35").Value = "  +7.59%  "
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").Formula = "'3.13"
$ws.Range("E38").Value = "  +11.35%  "
$ws.Range("E39").Value = "  +8.97%  "
$ws.Range("D40").Formula = "'2.10"
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Formula = "'9.09"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").Formula = "'0.309"
$ws.Range("E43").Value = "  +16.34%  "
$ws.Range("D44").Formula = "'44.19"
$ws.Range("E44").Value = "  +14.94%  "
$ws.Range("D45").Formula = "'392.11"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("D47").Value = "2.719.97"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Formula = "'132.73"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("D49").Formula = "'25.50"
$ws.Range("E49").Value = "  +10.95%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +8.36%  "
